$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the new row 193 D cell carries the same date number format as the rest of column D
$ws.Cells.Item(193, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 147
$ws.Cells.Item(147, 1).Value = 3
$ws.Cells.Item(147, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(147, 3).Value = "Coquimbo"
$ws.Cells.Item(147, 4).Value = 44463
$ws.Cells.Item(147, 5).Value = 5
$ws.Cells.Item(147, 6).Value = 100114013
$ws.Cells.Item(147, 7).Value = "Zanahoria"
$ws.Cells.Item(147, 8).Value = "Sin especificar"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 390
$ws.Cells.Item(147, 11).Value = 8500
$ws.Cells.Item(147, 12).Value = 9000
$ws.Cells.Item(147, 13).Value = 8705
$ws.Cells.Item(147, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(147, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(147, 16).Value = 435
$ws.Cells.Item(147, 17).Value = 20
$ws.Cells.Item(147, 18).Value = "Hortaliza"

# Row 148
$ws.Cells.Item(148, 1).Value = 3
$ws.Cells.Item(148, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(148, 3).Value = "Coquimbo"
$ws.Cells.Item(148, 4).Value = 44196
$ws.Cells.Item(148, 5).Value = 5
$ws.Cells.Item(148, 6).Value = 100114013
$ws.Cells.Item(148, 7).Value = "Zanahoria"
$ws.Cells.Item(148, 8).Value = "Sin especificar"
$ws.Cells.Item(148, 9).Value = "Primera"
$ws.Cells.Item(148, 10).Value = 190
$ws.Cells.Item(148, 11).Value = 6000
$ws.Cells.Item(148, 12).Value = 6000
$ws.Cells.Item(148, 13).Value = 6000
$ws.Cells.Item(148, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(148, 15).Value = "Chillán"
$ws.Cells.Item(148, 16).Value = 300
$ws.Cells.Item(148, 17).Value = 20
$ws.Cells.Item(148, 18).Value = "Hortaliza"

# Row 149
$ws.Cells.Item(149, 1).Value = 3
$ws.Cells.Item(149, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(149, 3).Value = "Coquimbo"
$ws.Cells.Item(149, 4).Value = 44301
$ws.Cells.Item(149, 5).Value = 5
$ws.Cells.Item(149, 6).Value = 100114013
$ws.Cells.Item(149, 7).Value = "Zanahoria"
$ws.Cells.Item(149, 8).Value = "Sin especificar"
$ws.Cells.Item(149, 9).Value = "Primera"
$ws.Cells.Item(149, 10).Value = 630
$ws.Cells.Item(149, 11).Value = 6000
$ws.Cells.Item(149, 12).Value = 6500
$ws.Cells.Item(149, 13).Value = 6222
$ws.Cells.Item(149, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(149, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(149, 16).Value = 311
$ws.Cells.Item(149, 17).Value = 20
$ws.Cells.Item(149, 18).Value = "Hortaliza"

# Row 150
$ws.Cells.Item(150, 1).Value = 3
$ws.Cells.Item(150, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(150, 3).Value = "Coquimbo"
$ws.Cells.Item(150, 4).Value = 44251
$ws.Cells.Item(150, 5).Value = 5
$ws.Cells.Item(150, 6).Value = 100114013
$ws.Cells.Item(150, 7).Value = "Zanahoria"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 180
$ws.Cells.Item(150, 11).Value = 6500
$ws.Cells.Item(150, 12).Value = 6500
$ws.Cells.Item(150, 13).Value = 6500
$ws.Cells.Item(150, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(150, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(150, 16).Value = 325
$ws.Cells.Item(150, 17).Value = 20
$ws.Cells.Item(150, 18).Value = "Hortaliza"

# Row 151
$ws.Cells.Item(151, 1).Value = 3
$ws.Cells.Item(151, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(151, 3).Value = "Coquimbo"
$ws.Cells.Item(151, 4).Value = 44243
$ws.Cells.Item(151, 5).Value = 5
$ws.Cells.Item(151, 6).Value = 100114013
$ws.Cells.Item(151, 7).Value = "Zanahoria"
$ws.Cells.Item(151, 8).Value = "Sin especificar"
$ws.Cells.Item(151, 9).Value = "Primera"
$ws.Cells.Item(151, 10).Value = 85
$ws.Cells.Item(151, 11).Value = 6500
$ws.Cells.Item(151, 12).Value = 6500
$ws.Cells.Item(151, 13).Value = 6500
$ws.Cells.Item(151, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(151, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(151, 16).Value = 325
$ws.Cells.Item(151, 17).Value = 20
$ws.Cells.Item(151, 18).Value = "Hortaliza"

# Row 152
$ws.Cells.Item(152, 1).Value = 3
$ws.Cells.Item(152, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(152, 3).Value = "Coquimbo"
$ws.Cells.Item(152, 4).Value = 44252
$ws.Cells.Item(152, 5).Value = 5
$ws.Cells.Item(152, 6).Value = 100114013
$ws.Cells.Item(152, 7).Value = "Zanahoria"
$ws.Cells.Item(152, 8).Value = "Sin especificar"
$ws.Cells.Item(152, 9).Value = "Primera"
$ws.Cells.Item(152, 10).Value = 160
$ws.Cells.Item(152, 11).Value = 6000
$ws.Cells.Item(152, 12).Value = 6000
$ws.Cells.Item(152, 13).Value = 6000
$ws.Cells.Item(152, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(152, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(152, 16).Value = 300
$ws.Cells.Item(152, 17).Value = 20
$ws.Cells.Item(152, 18).Value = "Hortaliza"

# Row 153
$ws.Cells.Item(153, 1).Value = 3
$ws.Cells.Item(153, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(153, 3).Value = "Coquimbo"
$ws.Cells.Item(153, 4).Value = 44166
$ws.Cells.Item(153, 5).Value = 5
$ws.Cells.Item(153, 6).Value = 100114013
$ws.Cells.Item(153, 7).Value = "Zanahoria"
$ws.Cells.Item(153, 8).Value = "Sin especificar"
$ws.Cells.Item(153, 9).Value = "Primera"
$ws.Cells.Item(153, 10).Value = 68
$ws.Cells.Item(153, 11).Value = 5000
$ws.Cells.Item(153, 12).Value = 5000
$ws.Cells.Item(153, 13).Value = 5000
$ws.Cells.Item(153, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(153, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(153, 16).Value = 250
$ws.Cells.Item(153, 17).Value = 20
$ws.Cells.Item(153, 18).Value = "Hortaliza"

# Row 154
$ws.Cells.Item(154, 1).Value = 3
$ws.Cells.Item(154, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(154, 3).Value = "Coquimbo"
$ws.Cells.Item(154, 4).Value = 44168
$ws.Cells.Item(154, 5).Value = 5
$ws.Cells.Item(154, 6).Value = 100114013
$ws.Cells.Item(154, 7).Value = "Zanahoria"
$ws.Cells.Item(154, 8).Value = "Sin especificar"
$ws.Cells.Item(154, 9).Value = "Primera"
$ws.Cells.Item(154, 10).Value = 170
$ws.Cells.Item(154, 11).Value = 5500
$ws.Cells.Item(154, 12).Value = 6000
$ws.Cells.Item(154, 13).Value = 5735
$ws.Cells.Item(154, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(154, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(154, 16).Value = 287
$ws.Cells.Item(154, 17).Value = 20
$ws.Cells.Item(154, 18).Value = "Hortaliza"

# Row 155
$ws.Cells.Item(155, 1).Value = 3
$ws.Cells.Item(155, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(155, 3).Value = "Coquimbo"
$ws.Cells.Item(155, 4).Value = 44369
$ws.Cells.Item(155, 5).Value = 5
$ws.Cells.Item(155, 6).Value = 100114013
$ws.Cells.Item(155, 7).Value = "Zanahoria"
$ws.Cells.Item(155, 8).Value = "Sin especificar"
$ws.Cells.Item(155, 9).Value = "Primera"
$ws.Cells.Item(155, 10).Value = 250
$ws.Cells.Item(155, 11).Value = 5500
$ws.Cells.Item(155, 12).Value = 5500
$ws.Cells.Item(155, 13).Value = 5500
$ws.Cells.Item(155, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(155, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(155, 16).Value = 275
$ws.Cells.Item(155, 17).Value = 20
$ws.Cells.Item(155, 18).Value = "Hortaliza"

# Row 156
$ws.Cells.Item(156, 1).Value = 3
$ws.Cells.Item(156, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(156, 3).Value = "Coquimbo"
$ws.Cells.Item(156, 4).Value = 44433
$ws.Cells.Item(156, 5).Value = 5
$ws.Cells.Item(156, 6).Value = 100114013
$ws.Cells.Item(156, 7).Value = "Zanahoria"
$ws.Cells.Item(156, 8).Value = "Sin especificar"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 310
$ws.Cells.Item(156, 11).Value = 5000
$ws.Cells.Item(156, 12).Value = 5500
$ws.Cells.Item(156, 13).Value = 5242
$ws.Cells.Item(156, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(156, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(156, 16).Value = 262
$ws.Cells.Item(156, 17).Value = 20
$ws.Cells.Item(156, 18).Value = "Hortaliza"

# Row 157
$ws.Cells.Item(157, 1).Value = 3
$ws.Cells.Item(157, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(157, 3).Value = "Coquimbo"
$ws.Cells.Item(157, 4).Value = 44221
$ws.Cells.Item(157, 5).Value = 5
$ws.Cells.Item(157, 6).Value = 100114013
$ws.Cells.Item(157, 7).Value = "Zanahoria"
$ws.Cells.Item(157, 8).Value = "Sin especificar"
$ws.Cells.Item(157, 9).Value = "Primera"
$ws.Cells.Item(157, 10).Value = 200
$ws.Cells.Item(157, 11).Value = 6500
$ws.Cells.Item(157, 12).Value = 6500
$ws.Cells.Item(157, 13).Value = 6500
$ws.Cells.Item(157, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(157, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(157, 16).Value = 325
$ws.Cells.Item(157, 17).Value = 20
$ws.Cells.Item(157, 18).Value = "Hortaliza"

# Row 158
$ws.Cells.Item(158, 1).Value = 3
$ws.Cells.Item(158, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(158, 3).Value = "Coquimbo"
$ws.Cells.Item(158, 4).Value = 44371
$ws.Cells.Item(158, 5).Value = 5
$ws.Cells.Item(158, 6).Value = 100114013
$ws.Cells.Item(158, 7).Value = "Zanahoria"
$ws.Cells.Item(158, 8).Value = "Sin especificar"
$ws.Cells.Item(158, 9).Value = "Primera"
$ws.Cells.Item(158, 10).Value = 430
$ws.Cells.Item(158, 11).Value = 5000
$ws.Cells.Item(158, 12).Value = 5500
$ws.Cells.Item(158, 13).Value = 5291
$ws.Cells.Item(158, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(158, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(158, 16).Value = 265
$ws.Cells.Item(158, 17).Value = 20
$ws.Cells.Item(158, 18).Value = "Hortaliza"

# Row 159
$ws.Cells.Item(159, 1).Value = 3
$ws.Cells.Item(159, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(159, 3).Value = "Coquimbo"
$ws.Cells.Item(159, 4).Value = 44316
$ws.Cells.Item(159, 5).Value = 5
$ws.Cells.Item(159, 6).Value = 100114013
$ws.Cells.Item(159, 7).Value = "Zanahoria"
$ws.Cells.Item(159, 8).Value = "Sin especificar"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 270
$ws.Cells.Item(159, 11).Value = 5500
$ws.Cells.Item(159, 12).Value = 6000
$ws.Cells.Item(159, 13).Value = 5722
$ws.Cells.Item(159, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(159, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(159, 16).Value = 286
$ws.Cells.Item(159, 17).Value = 20
$ws.Cells.Item(159, 18).Value = "Hortaliza"

# Row 160
$ws.Cells.Item(160, 1).Value = 3
$ws.Cells.Item(160, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(160, 3).Value = "Coquimbo"
$ws.Cells.Item(160, 4).Value = 44279
$ws.Cells.Item(160, 5).Value = 5
$ws.Cells.Item(160, 6).Value = 100114013
$ws.Cells.Item(160, 7).Value = "Zanahoria"
$ws.Cells.Item(160, 8).Value = "Sin especificar"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 450
$ws.Cells.Item(160, 11).Value = 6000
$ws.Cells.Item(160, 12).Value = 6500
$ws.Cells.Item(160, 13).Value = 6222
$ws.Cells.Item(160, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(160, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(160, 16).Value = 311
$ws.Cells.Item(160, 17).Value = 20
$ws.Cells.Item(160, 18).Value = "Hortaliza"

# Row 161
$ws.Cells.Item(161, 1).Value = 3
$ws.Cells.Item(161, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(161, 3).Value = "Coquimbo"
$ws.Cells.Item(161, 4).Value = 44397
$ws.Cells.Item(161, 5).Value = 5
$ws.Cells.Item(161, 6).Value = 100114013
$ws.Cells.Item(161, 7).Value = "Zanahoria"
$ws.Cells.Item(161, 8).Value = "Sin especificar"
$ws.Cells.Item(161, 9).Value = "Primera"
$ws.Cells.Item(161, 10).Value = 180
$ws.Cells.Item(161, 11).Value = 5500
$ws.Cells.Item(161, 12).Value = 5500
$ws.Cells.Item(161, 13).Value = 5500
$ws.Cells.Item(161, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(161, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(161, 16).Value = 275
$ws.Cells.Item(161, 17).Value = 20
$ws.Cells.Item(161, 18).Value = "Hortaliza"

# Row 162
$ws.Cells.Item(162, 1).Value = 3
$ws.Cells.Item(162, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(162, 3).Value = "Coquimbo"
$ws.Cells.Item(162, 4).Value = 44363
$ws.Cells.Item(162, 5).Value = 5
$ws.Cells.Item(162, 6).Value = 100114013
$ws.Cells.Item(162, 7).Value = "Zanahoria"
$ws.Cells.Item(162, 8).Value = "Sin especificar"
$ws.Cells.Item(162, 9).Value = "Primera"
$ws.Cells.Item(162, 10).Value = 340
$ws.Cells.Item(162, 11).Value = 5500
$ws.Cells.Item(162, 12).Value = 6000
$ws.Cells.Item(162, 13).Value = 5765
$ws.Cells.Item(162, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(162, 15).Value = "Chillán"
$ws.Cells.Item(162, 16).Value = 288
$ws.Cells.Item(162, 17).Value = 20
$ws.Cells.Item(162, 18).Value = "Hortaliza"

# Row 163
$ws.Cells.Item(163, 1).Value = 3
$ws.Cells.Item(163, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(163, 3).Value = "Coquimbo"
$ws.Cells.Item(163, 4).Value = 44277
$ws.Cells.Item(163, 5).Value = 5
$ws.Cells.Item(163, 6).Value = 100114013
$ws.Cells.Item(163, 7).Value = "Zanahoria"
$ws.Cells.Item(163, 8).Value = "Sin especificar"
$ws.Cells.Item(163, 9).Value = "Primera"
$ws.Cells.Item(163, 10).Value = 310
$ws.Cells.Item(163, 11).Value = 5500
$ws.Cells.Item(163, 12).Value = 6000
$ws.Cells.Item(163, 13).Value = 5758
$ws.Cells.Item(163, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(163, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(163, 16).Value = 288
$ws.Cells.Item(163, 17).Value = 20
$ws.Cells.Item(163, 18).Value = "Hortaliza"

# Row 164
$ws.Cells.Item(164, 1).Value = 3
$ws.Cells.Item(164, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(164, 3).Value = "Coquimbo"
$ws.Cells.Item(164, 4).Value = 44291
$ws.Cells.Item(164, 5).Value = 5
$ws.Cells.Item(164, 6).Value = 100114013
$ws.Cells.Item(164, 7).Value = "Zanahoria"
$ws.Cells.Item(164, 8).Value = "Sin especificar"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 220
$ws.Cells.Item(164, 11).Value = 6000
$ws.Cells.Item(164, 12).Value = 6000
$ws.Cells.Item(164, 13).Value = 6000
$ws.Cells.Item(164, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(164, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(164, 16).Value = 300
$ws.Cells.Item(164, 17).Value = 20
$ws.Cells.Item(164, 18).Value = "Hortaliza"

# Row 165
$ws.Cells.Item(165, 1).Value = 3
$ws.Cells.Item(165, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(165, 3).Value = "Coquimbo"
$ws.Cells.Item(165, 4).Value = 44273
$ws.Cells.Item(165, 5).Value = 5
$ws.Cells.Item(165, 6).Value = 100114013
$ws.Cells.Item(165, 7).Value = "Zanahoria"
$ws.Cells.Item(165, 8).Value = "Sin especificar"
$ws.Cells.Item(165, 9).Value = "Primera"
$ws.Cells.Item(165, 10).Value = 340
$ws.Cells.Item(165, 11).Value = 6000
$ws.Cells.Item(165, 12).Value = 6500
$ws.Cells.Item(165, 13).Value = 6265
$ws.Cells.Item(165, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(165, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(165, 16).Value = 313
$ws.Cells.Item(165, 17).Value = 20
$ws.Cells.Item(165, 18).Value = "Hortaliza"

# Row 166
$ws.Cells.Item(166, 1).Value = 3
$ws.Cells.Item(166, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(166, 3).Value = "Coquimbo"
$ws.Cells.Item(166, 4).Value = 44438
$ws.Cells.Item(166, 5).Value = 5
$ws.Cells.Item(166, 6).Value = 100114013
$ws.Cells.Item(166, 7).Value = "Zanahoria"
$ws.Cells.Item(166, 8).Value = "Sin especificar"
$ws.Cells.Item(166, 9).Value = "Primera"
$ws.Cells.Item(166, 10).Value = 440
$ws.Cells.Item(166, 11).Value = 5000
$ws.Cells.Item(166, 12).Value = 5500
$ws.Cells.Item(166, 13).Value = 5205
$ws.Cells.Item(166, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(166, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(166, 16).Value = 260
$ws.Cells.Item(166, 17).Value = 20
$ws.Cells.Item(166, 18).Value = "Hortaliza"

# Row 167
$ws.Cells.Item(167, 1).Value = 3
$ws.Cells.Item(167, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(167, 3).Value = "Coquimbo"
$ws.Cells.Item(167, 4).Value = 44372
$ws.Cells.Item(167, 5).Value = 5
$ws.Cells.Item(167, 6).Value = 100114013
$ws.Cells.Item(167, 7).Value = "Zanahoria"
$ws.Cells.Item(167, 8).Value = "Sin especificar"
$ws.Cells.Item(167, 9).Value = "Primera"
$ws.Cells.Item(167, 10).Value = 280
$ws.Cells.Item(167, 11).Value = 5500
$ws.Cells.Item(167, 12).Value = 5500
$ws.Cells.Item(167, 13).Value = 5500
$ws.Cells.Item(167, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(167, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(167, 16).Value = 275
$ws.Cells.Item(167, 17).Value = 20
$ws.Cells.Item(167, 18).Value = "Hortaliza"

# Row 168
$ws.Cells.Item(168, 1).Value = 3
$ws.Cells.Item(168, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(168, 3).Value = "Coquimbo"
$ws.Cells.Item(168, 4).Value = 44286
$ws.Cells.Item(168, 5).Value = 5
$ws.Cells.Item(168, 6).Value = 100114013
$ws.Cells.Item(168, 7).Value = "Zanahoria"
$ws.Cells.Item(168, 8).Value = "Sin especificar"
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 280
$ws.Cells.Item(168, 11).Value = 6000
$ws.Cells.Item(168, 12).Value = 6000
$ws.Cells.Item(168, 13).Value = 6000
$ws.Cells.Item(168, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(168, 15).Value = "Chillán"
$ws.Cells.Item(168, 16).Value = 300
$ws.Cells.Item(168, 17).Value = 20
$ws.Cells.Item(168, 18).Value = "Hortaliza"

# Row 169
$ws.Cells.Item(169, 1).Value = 3
$ws.Cells.Item(169, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(169, 3).Value = "Coquimbo"
$ws.Cells.Item(169, 4).Value = 44209
$ws.Cells.Item(169, 5).Value = 5
$ws.Cells.Item(169, 6).Value = 100114013
$ws.Cells.Item(169, 7).Value = "Zanahoria"
$ws.Cells.Item(169, 8).Value = "Sin especificar"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 180
$ws.Cells.Item(169, 11).Value = 6500
$ws.Cells.Item(169, 12).Value = 6500
$ws.Cells.Item(169, 13).Value = 6500
$ws.Cells.Item(169, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(169, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(169, 16).Value = 325
$ws.Cells.Item(169, 17).Value = 20
$ws.Cells.Item(169, 18).Value = "Hortaliza"

# Row 170
$ws.Cells.Item(170, 1).Value = 3
$ws.Cells.Item(170, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(170, 3).Value = "Coquimbo"
$ws.Cells.Item(170, 4).Value = 44356
$ws.Cells.Item(170, 5).Value = 5
$ws.Cells.Item(170, 6).Value = 100114013
$ws.Cells.Item(170, 7).Value = "Zanahoria"
$ws.Cells.Item(170, 8).Value = "Sin especificar"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 320
$ws.Cells.Item(170, 11).Value = 5000
$ws.Cells.Item(170, 12).Value = 5500
$ws.Cells.Item(170, 13).Value = 5250
$ws.Cells.Item(170, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(170, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(170, 16).Value = 262
$ws.Cells.Item(170, 17).Value = 20
$ws.Cells.Item(170, 18).Value = "Hortaliza"

# Row 171
$ws.Cells.Item(171, 1).Value = 3
$ws.Cells.Item(171, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(171, 3).Value = "Coquimbo"
$ws.Cells.Item(171, 4).Value = 44160
$ws.Cells.Item(171, 5).Value = 5
$ws.Cells.Item(171, 6).Value = 100114013
$ws.Cells.Item(171, 7).Value = "Zanahoria"
$ws.Cells.Item(171, 8).Value = "Sin especificar"
$ws.Cells.Item(171, 9).Value = "Primera"
$ws.Cells.Item(171, 10).Value = 160
$ws.Cells.Item(171, 11).Value = 5000
$ws.Cells.Item(171, 12).Value = 5000
$ws.Cells.Item(171, 13).Value = 5000
$ws.Cells.Item(171, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(171, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(171, 16).Value = 250
$ws.Cells.Item(171, 17).Value = 20
$ws.Cells.Item(171, 18).Value = "Hortaliza"

# Row 172
$ws.Cells.Item(172, 1).Value = 3
$ws.Cells.Item(172, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(172, 3).Value = "Coquimbo"
$ws.Cells.Item(172, 4).Value = 44351
$ws.Cells.Item(172, 5).Value = 5
$ws.Cells.Item(172, 6).Value = 100114013
$ws.Cells.Item(172, 7).Value = "Zanahoria"
$ws.Cells.Item(172, 8).Value = "Sin especificar"
$ws.Cells.Item(172, 9).Value = "Primera"
$ws.Cells.Item(172, 10).Value = 510
$ws.Cells.Item(172, 11).Value = 5000
$ws.Cells.Item(172, 12).Value = 5500
$ws.Cells.Item(172, 13).Value = 5245
$ws.Cells.Item(172, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(172, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(172, 16).Value = 262
$ws.Cells.Item(172, 17).Value = 20
$ws.Cells.Item(172, 18).Value = "Hortaliza"

# Row 173
$ws.Cells.Item(173, 1).Value = 3
$ws.Cells.Item(173, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(173, 3).Value = "Coquimbo"
$ws.Cells.Item(173, 4).Value = 44365
$ws.Cells.Item(173, 5).Value = 5
$ws.Cells.Item(173, 6).Value = 100114013
$ws.Cells.Item(173, 7).Value = "Zanahoria"
$ws.Cells.Item(173, 8).Value = "Sin especificar"
$ws.Cells.Item(173, 9).Value = "Primera"
$ws.Cells.Item(173, 10).Value = 450
$ws.Cells.Item(173, 11).Value = 5500
$ws.Cells.Item(173, 12).Value = 6000
$ws.Cells.Item(173, 13).Value = 5756
$ws.Cells.Item(173, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(173, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(173, 16).Value = 288
$ws.Cells.Item(173, 17).Value = 20
$ws.Cells.Item(173, 18).Value = "Hortaliza"

# Row 174
$ws.Cells.Item(174, 1).Value = 3
$ws.Cells.Item(174, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(174, 3).Value = "Coquimbo"
$ws.Cells.Item(174, 4).Value = 44306
$ws.Cells.Item(174, 5).Value = 5
$ws.Cells.Item(174, 6).Value = 100114013
$ws.Cells.Item(174, 7).Value = "Zanahoria"
$ws.Cells.Item(174, 8).Value = "Sin especificar"
$ws.Cells.Item(174, 9).Value = "Primera"
$ws.Cells.Item(174, 10).Value = 160
$ws.Cells.Item(174, 11).Value = 5500
$ws.Cells.Item(174, 12).Value = 5500
$ws.Cells.Item(174, 13).Value = 5500
$ws.Cells.Item(174, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(174, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(174, 16).Value = 275
$ws.Cells.Item(174, 17).Value = 20
$ws.Cells.Item(174, 18).Value = "Hortaliza"

# Row 175
$ws.Cells.Item(175, 1).Value = 3
$ws.Cells.Item(175, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(175, 3).Value = "Coquimbo"
$ws.Cells.Item(175, 4).Value = 44215
$ws.Cells.Item(175, 5).Value = 5
$ws.Cells.Item(175, 6).Value = 100114013
$ws.Cells.Item(175, 7).Value = "Zanahoria"
$ws.Cells.Item(175, 8).Value = "Sin especificar"
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 160
$ws.Cells.Item(175, 11).Value = 7000
$ws.Cells.Item(175, 12).Value = 7000
$ws.Cells.Item(175, 13).Value = 7000
$ws.Cells.Item(175, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(175, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(175, 16).Value = 350
$ws.Cells.Item(175, 17).Value = 20
$ws.Cells.Item(175, 18).Value = "Hortaliza"

# Row 176
$ws.Cells.Item(176, 1).Value = 3
$ws.Cells.Item(176, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(176, 3).Value = "Coquimbo"
$ws.Cells.Item(176, 4).Value = 44175
$ws.Cells.Item(176, 5).Value = 5
$ws.Cells.Item(176, 6).Value = 100114013
$ws.Cells.Item(176, 7).Value = "Zanahoria"
$ws.Cells.Item(176, 8).Value = "Sin especificar"
$ws.Cells.Item(176, 9).Value = "Primera"
$ws.Cells.Item(176, 10).Value = 160
$ws.Cells.Item(176, 11).Value = 5000
$ws.Cells.Item(176, 12).Value = 5000
$ws.Cells.Item(176, 13).Value = 5000
$ws.Cells.Item(176, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(176, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(176, 16).Value = 250
$ws.Cells.Item(176, 17).Value = 20
$ws.Cells.Item(176, 18).Value = "Hortaliza"

# Row 177
$ws.Cells.Item(177, 1).Value = 3
$ws.Cells.Item(177, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(177, 3).Value = "Coquimbo"
$ws.Cells.Item(177, 4).Value = 44175
$ws.Cells.Item(177, 5).Value = 5
$ws.Cells.Item(177, 6).Value = 100114013
$ws.Cells.Item(177, 7).Value = "Zanahoria"
$ws.Cells.Item(177, 8).Value = "Sin especificar"
$ws.Cells.Item(177, 9).Value = "Segunda"
$ws.Cells.Item(177, 10).Value = 180
$ws.Cells.Item(177, 11).Value = 4000
$ws.Cells.Item(177, 12).Value = 4000
$ws.Cells.Item(177, 13).Value = 4000
$ws.Cells.Item(177, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(177, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(177, 16).Value = 200
$ws.Cells.Item(177, 17).Value = 20
$ws.Cells.Item(177, 18).Value = "Hortaliza"

# Row 178
$ws.Cells.Item(178, 1).Value = 3
$ws.Cells.Item(178, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(178, 3).Value = "Coquimbo"
$ws.Cells.Item(178, 4).Value = 44357
$ws.Cells.Item(178, 5).Value = 5
$ws.Cells.Item(178, 6).Value = 100114013
$ws.Cells.Item(178, 7).Value = "Zanahoria"
$ws.Cells.Item(178, 8).Value = "Sin especificar"
$ws.Cells.Item(178, 9).Value = "Primera"
$ws.Cells.Item(178, 10).Value = 310
$ws.Cells.Item(178, 11).Value = 5500
$ws.Cells.Item(178, 12).Value = 5800
$ws.Cells.Item(178, 13).Value = 5645
$ws.Cells.Item(178, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(178, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(178, 16).Value = 282
$ws.Cells.Item(178, 17).Value = 20
$ws.Cells.Item(178, 18).Value = "Hortaliza"

# Row 179
$ws.Cells.Item(179, 1).Value = 3
$ws.Cells.Item(179, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(179, 3).Value = "Coquimbo"
$ws.Cells.Item(179, 4).Value = 44203
$ws.Cells.Item(179, 5).Value = 5
$ws.Cells.Item(179, 6).Value = 100114013
$ws.Cells.Item(179, 7).Value = "Zanahoria"
$ws.Cells.Item(179, 8).Value = "Sin especificar"
$ws.Cells.Item(179, 9).Value = "Primera"
$ws.Cells.Item(179, 10).Value = 220
$ws.Cells.Item(179, 11).Value = 6000
$ws.Cells.Item(179, 12).Value = 6000
$ws.Cells.Item(179, 13).Value = 6000
$ws.Cells.Item(179, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(179, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(179, 16).Value = 300
$ws.Cells.Item(179, 17).Value = 20
$ws.Cells.Item(179, 18).Value = "Hortaliza"

# Row 180
$ws.Cells.Item(180, 1).Value = 3
$ws.Cells.Item(180, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(180, 3).Value = "Coquimbo"
$ws.Cells.Item(180, 4).Value = 44162
$ws.Cells.Item(180, 5).Value = 5
$ws.Cells.Item(180, 6).Value = 100114013
$ws.Cells.Item(180, 7).Value = "Zanahoria"
$ws.Cells.Item(180, 8).Value = "Sin especificar"
$ws.Cells.Item(180, 9).Value = "Primera"
$ws.Cells.Item(180, 10).Value = 130
$ws.Cells.Item(180, 11).Value = 5000
$ws.Cells.Item(180, 12).Value = 5000
$ws.Cells.Item(180, 13).Value = 5000
$ws.Cells.Item(180, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(180, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(180, 16).Value = 250
$ws.Cells.Item(180, 17).Value = 20
$ws.Cells.Item(180, 18).Value = "Hortaliza"

# Row 181
$ws.Cells.Item(181, 1).Value = 3
$ws.Cells.Item(181, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(181, 3).Value = "Coquimbo"
$ws.Cells.Item(181, 4).Value = 44410
$ws.Cells.Item(181, 5).Value = 5
$ws.Cells.Item(181, 6).Value = 100114013
$ws.Cells.Item(181, 7).Value = "Zanahoria"
$ws.Cells.Item(181, 8).Value = "Sin especificar"
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 510
$ws.Cells.Item(181, 11).Value = 4500
$ws.Cells.Item(181, 12).Value = 5000
$ws.Cells.Item(181, 13).Value = 4745
$ws.Cells.Item(181, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(181, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(181, 16).Value = 237
$ws.Cells.Item(181, 17).Value = 20
$ws.Cells.Item(181, 18).Value = "Hortaliza"

# Row 182
$ws.Cells.Item(182, 1).Value = 3
$ws.Cells.Item(182, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(182, 3).Value = "Coquimbo"
$ws.Cells.Item(182, 4).Value = 44411
$ws.Cells.Item(182, 5).Value = 5
$ws.Cells.Item(182, 6).Value = 100114013
$ws.Cells.Item(182, 7).Value = "Zanahoria"
$ws.Cells.Item(182, 8).Value = "Sin especificar"
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 410
$ws.Cells.Item(182, 11).Value = 4500
$ws.Cells.Item(182, 12).Value = 5000
$ws.Cells.Item(182, 13).Value = 4805
$ws.Cells.Item(182, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(182, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(182, 16).Value = 240
$ws.Cells.Item(182, 17).Value = 20
$ws.Cells.Item(182, 18).Value = "Hortaliza"

# Row 183
$ws.Cells.Item(183, 1).Value = 3
$ws.Cells.Item(183, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(183, 3).Value = "Coquimbo"
$ws.Cells.Item(183, 4).Value = 44257
$ws.Cells.Item(183, 5).Value = 5
$ws.Cells.Item(183, 6).Value = 100114013
$ws.Cells.Item(183, 7).Value = "Zanahoria"
$ws.Cells.Item(183, 8).Value = "Sin especificar"
$ws.Cells.Item(183, 9).Value = "Primera"
$ws.Cells.Item(183, 10).Value = 280
$ws.Cells.Item(183, 11).Value = 6000
$ws.Cells.Item(183, 12).Value = 6000
$ws.Cells.Item(183, 13).Value = 6000
$ws.Cells.Item(183, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(183, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(183, 16).Value = 300
$ws.Cells.Item(183, 17).Value = 20
$ws.Cells.Item(183, 18).Value = "Hortaliza"

# Row 184
$ws.Cells.Item(184, 1).Value = 3
$ws.Cells.Item(184, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(184, 3).Value = "Coquimbo"
$ws.Cells.Item(184, 4).Value = 44244
$ws.Cells.Item(184, 5).Value = 5
$ws.Cells.Item(184, 6).Value = 100114013
$ws.Cells.Item(184, 7).Value = "Zanahoria"
$ws.Cells.Item(184, 8).Value = "Sin especificar"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 310
$ws.Cells.Item(184, 11).Value = 6000
$ws.Cells.Item(184, 12).Value = 6500
$ws.Cells.Item(184, 13).Value = 6242
$ws.Cells.Item(184, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(184, 15).Value = "Chillán"
$ws.Cells.Item(184, 16).Value = 312
$ws.Cells.Item(184, 17).Value = 20
$ws.Cells.Item(184, 18).Value = "Hortaliza"

# Row 185
$ws.Cells.Item(185, 1).Value = 3
$ws.Cells.Item(185, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(185, 3).Value = "Coquimbo"
$ws.Cells.Item(185, 4).Value = 44176
$ws.Cells.Item(185, 5).Value = 5
$ws.Cells.Item(185, 6).Value = 100114013
$ws.Cells.Item(185, 7).Value = "Zanahoria"
$ws.Cells.Item(185, 8).Value = "Sin especificar"
$ws.Cells.Item(185, 9).Value = "Primera"
$ws.Cells.Item(185, 10).Value = 310
$ws.Cells.Item(185, 11).Value = 5000
$ws.Cells.Item(185, 12).Value = 5500
$ws.Cells.Item(185, 13).Value = 5258
$ws.Cells.Item(185, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(185, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(185, 16).Value = 263
$ws.Cells.Item(185, 17).Value = 20
$ws.Cells.Item(185, 18).Value = "Hortaliza"

# Row 186
$ws.Cells.Item(186, 1).Value = 3
$ws.Cells.Item(186, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(186, 3).Value = "Coquimbo"
$ws.Cells.Item(186, 4).Value = 44176
$ws.Cells.Item(186, 5).Value = 5
$ws.Cells.Item(186, 6).Value = 100114013
$ws.Cells.Item(186, 7).Value = "Zanahoria"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "Segunda"
$ws.Cells.Item(186, 10).Value = 50
$ws.Cells.Item(186, 11).Value = 4000
$ws.Cells.Item(186, 12).Value = 4000
$ws.Cells.Item(186, 13).Value = 4000
$ws.Cells.Item(186, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(186, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(186, 16).Value = 200
$ws.Cells.Item(186, 17).Value = 20
$ws.Cells.Item(186, 18).Value = "Hortaliza"

# Row 187
$ws.Cells.Item(187, 1).Value = 3
$ws.Cells.Item(187, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(187, 3).Value = "Coquimbo"
$ws.Cells.Item(187, 4).Value = 44239
$ws.Cells.Item(187, 5).Value = 5
$ws.Cells.Item(187, 6).Value = 100114013
$ws.Cells.Item(187, 7).Value = "Zanahoria"
$ws.Cells.Item(187, 8).Value = "Sin especificar"
$ws.Cells.Item(187, 9).Value = "Primera"
$ws.Cells.Item(187, 10).Value = 240
$ws.Cells.Item(187, 11).Value = 6500
$ws.Cells.Item(187, 12).Value = 7000
$ws.Cells.Item(187, 13).Value = 6667
$ws.Cells.Item(187, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(187, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(187, 16).Value = 333
$ws.Cells.Item(187, 17).Value = 20
$ws.Cells.Item(187, 18).Value = "Hortaliza"

# Row 188
$ws.Cells.Item(188, 1).Value = 3
$ws.Cells.Item(188, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(188, 3).Value = "Coquimbo"
$ws.Cells.Item(188, 4).Value = 44376
$ws.Cells.Item(188, 5).Value = 5
$ws.Cells.Item(188, 6).Value = 100114013
$ws.Cells.Item(188, 7).Value = "Zanahoria"
$ws.Cells.Item(188, 8).Value = "Sin especificar"
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 530
$ws.Cells.Item(188, 11).Value = 5500
$ws.Cells.Item(188, 12).Value = 5800
$ws.Cells.Item(188, 13).Value = 5642
$ws.Cells.Item(188, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(188, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(188, 16).Value = 282
$ws.Cells.Item(188, 17).Value = 20
$ws.Cells.Item(188, 18).Value = "Hortaliza"

# Row 189
$ws.Cells.Item(189, 1).Value = 3
$ws.Cells.Item(189, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(189, 3).Value = "Coquimbo"
$ws.Cells.Item(189, 4).Value = 44292
$ws.Cells.Item(189, 5).Value = 5
$ws.Cells.Item(189, 6).Value = 100114013
$ws.Cells.Item(189, 7).Value = "Zanahoria"
$ws.Cells.Item(189, 8).Value = "Sin especificar"
$ws.Cells.Item(189, 9).Value = "Primera"
$ws.Cells.Item(189, 10).Value = 85
$ws.Cells.Item(189, 11).Value = 6000
$ws.Cells.Item(189, 12).Value = 6000
$ws.Cells.Item(189, 13).Value = 6000
$ws.Cells.Item(189, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(189, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(189, 16).Value = 300
$ws.Cells.Item(189, 17).Value = 20
$ws.Cells.Item(189, 18).Value = "Hortaliza"

# Row 190
$ws.Cells.Item(190, 1).Value = 3
$ws.Cells.Item(190, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(190, 3).Value = "Coquimbo"
$ws.Cells.Item(190, 4).Value = 44358
$ws.Cells.Item(190, 5).Value = 5
$ws.Cells.Item(190, 6).Value = 100114013
$ws.Cells.Item(190, 7).Value = "Zanahoria"
$ws.Cells.Item(190, 8).Value = "Sin especificar"
$ws.Cells.Item(190, 9).Value = "Primera"
$ws.Cells.Item(190, 10).Value = 340
$ws.Cells.Item(190, 11).Value = 5500
$ws.Cells.Item(190, 12).Value = 6000
$ws.Cells.Item(190, 13).Value = 5765
$ws.Cells.Item(190, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(190, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(190, 16).Value = 288
$ws.Cells.Item(190, 17).Value = 20
$ws.Cells.Item(190, 18).Value = "Hortaliza"

# Row 191
$ws.Cells.Item(191, 1).Value = 3
$ws.Cells.Item(191, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(191, 3).Value = "Coquimbo"
$ws.Cells.Item(191, 4).Value = 44211
$ws.Cells.Item(191, 5).Value = 5
$ws.Cells.Item(191, 6).Value = 100114013
$ws.Cells.Item(191, 7).Value = "Zanahoria"
$ws.Cells.Item(191, 8).Value = "Sin especificar"
$ws.Cells.Item(191, 9).Value = "Primera"
$ws.Cells.Item(191, 10).Value = 160
$ws.Cells.Item(191, 11).Value = 8000
$ws.Cells.Item(191, 12).Value = 8000
$ws.Cells.Item(191, 13).Value = 8000
$ws.Cells.Item(191, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(191, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(191, 16).Value = 400
$ws.Cells.Item(191, 17).Value = 20
$ws.Cells.Item(191, 18).Value = "Hortaliza"

# Row 192
$ws.Cells.Item(192, 1).Value = 3
$ws.Cells.Item(192, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(192, 3).Value = "Coquimbo"
$ws.Cells.Item(192, 4).Value = 44425
$ws.Cells.Item(192, 5).Value = 5
$ws.Cells.Item(192, 6).Value = 100114013
$ws.Cells.Item(192, 7).Value = "Zanahoria"
$ws.Cells.Item(192, 8).Value = "Sin especificar"
$ws.Cells.Item(192, 9).Value = "Primera"
$ws.Cells.Item(192, 10).Value = 410
$ws.Cells.Item(192, 11).Value = 4500
$ws.Cells.Item(192, 12).Value = 5000
$ws.Cells.Item(192, 13).Value = 4780
$ws.Cells.Item(192, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(192, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(192, 16).Value = 239
$ws.Cells.Item(192, 17).Value = 20
$ws.Cells.Item(192, 18).Value = "Hortaliza"

# Row 193
$ws.Cells.Item(193, 1).Value = 3
$ws.Cells.Item(193, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(193, 3).Value = "Coquimbo"
$ws.Cells.Item(193, 4).Value = 44323
$ws.Cells.Item(193, 5).Value = 5
$ws.Cells.Item(193, 6).Value = 100114013
$ws.Cells.Item(193, 7).Value = "Zanahoria"
$ws.Cells.Item(193, 8).Value = "Sin especificar"
$ws.Cells.Item(193, 9).Value = "Primera"
$ws.Cells.Item(193, 10).Value = 300
$ws.Cells.Item(193, 11).Value = 5000
$ws.Cells.Item(193, 12).Value = 5500
$ws.Cells.Item(193, 13).Value = 5233
$ws.Cells.Item(193, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(193, 15).Value = "Chillán"
$ws.Cells.Item(193, 16).Value = 262
$ws.Cells.Item(193, 17).Value = 20
$ws.Cells.Item(193, 18).Value = "Hortaliza"
